{"js": "// Helper: wrap a <w:body> fragment in the full OOXML \"flat package\" envelope\n// that Range.insertOoxml / body.insertOoxml expect.\nfunction wrapOoxml(bodyXml) {\n  return '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n      '<pkg:part pkg:name=\"/word/document.xml\" ' +\n        'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n          '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n            '<w:body>' + bodyXml + '</w:body>' +\n          '</w:document>' +\n        '</pkg:xmlData>' +\n      '</pkg:part>' +\n    '</pkg:package>';\n}\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// --- Paragraph 1: \"A simple / demonstration / of a / query / :\" was split\n// across 5 runs (and preceded by leftover spell-check proofErr markers).\n// Collapse it back down to one clean run with the full sentence, keep the\n// existing bookmark, and drop the stale proofErr marks.\nconst firstRange = paragraphs.items[0].getRange();\nfirstRange.insertOoxml(\n  wrapOoxml(\n    '<w:p>' +\n      '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n      '<w:bookmarkEnd w:id=\"0\"/>' +\n      '<w:r><w:t>A simple demonstration of a query :</w:t></w:r>' +\n    '</w:p>'\n  ),\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// Reload the paragraph collection since paragraph 1 was structurally replaced.\nparagraphs.load(\"items\");\nawait context.sync();\n\n// --- Paragraph 2: \"anydsl\" loses its orange theme-color character formatting.\nconst secondRange = paragraphs.items[1].getRange();\nsecondRange.insertOoxml(\n  wrapOoxml('<w:p><w:r><w:t>anydsl</w:t></w:r></w:p>'),\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\nparagraphs.load(\"items\");\nawait context.sync();\n\n// --- Paragraph 4: the trailing empty paragraph had a stray empty <w:t/>\n// inside its run; normalize it down to a plain empty run.\nconst lastRange = paragraphs.items[3].getRange();\nlastRange.insertOoxml(\n  wrapOoxml('<w:p><w:r/></w:p>'),\n  Word.InsertLocation.replace\n);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Helper: wrap a <w:body> fragment in the full OOXML \"flat package\" envelope\n# that Range.InsertXML expects.\nfunction Wrap-Ooxml([string]$bodyXml) {\n    return @\"\n<?xml version=\"1.0\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>$bodyXml</w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n\"@\n}\n\n# --- Paragraph 1: \"A simple / demonstration / of a / query / :\" was split\n# across 5 runs (and preceded by leftover spell-check proofErr markers).\n# Collapse it back down to one clean run with the full sentence, keep the\n# existing bookmark, and drop the stale proofErr marks.\n$p1 = $d.Paragraphs.Item(1)\n$p1.Range.InsertXML((Wrap-Ooxml '<w:p><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/><w:r><w:t>A simple demonstration of a query :</w:t></w:r></w:p>')) | Out-Null\n\n# --- Paragraph 2: \"anydsl\" loses its orange theme-color character formatting.\n$p2 = $d.Paragraphs.Item(2)\n$p2.Range.InsertXML((Wrap-Ooxml '<w:p><w:r><w:t>anydsl</w:t></w:r></w:p>')) | Out-Null\n\n# --- Paragraph 4: the trailing empty paragraph had a stray empty <w:t/>\n# inside its run. It is also the document's very last paragraph, so its\n# Range only covers the paragraph mark itself; replacing that whole range\n# with a full <w:p> would duplicate the mark, so instead collapse to a\n# zero-length range just before the mark and insert a plain empty run.\n$p4 = $d.Paragraphs.Item(4)\n$r4 = $p4.Range\n$r4.MoveEnd(1, -1) | Out-Null\n$r4.InsertXML((Wrap-Ooxml '<w:r/>')) | Out-Null\n"}
